# "Lots of code cleanup"
# - Rebuild the second (mirrored) mini-table on Sheet1 so that it no longer
#   repeats the Date column (which duplicated column A); the 4 data columns
#   shift left from I:L to G:J.
# - Reset custom row heights on rows 24-26 and drop now-empty spacer rows
#   28/30/31.
# - Reflow/resize the chart, switch the page to portrait, and set a print
#   area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Wipe out the old H:L block (dates mirror + 4 data columns) so that
#    both values and styles (e.g. the date number format) are gone.
# ---------------------------------------------------------------------
$ws.Range("H1:L31").Clear()

# ---------------------------------------------------------------------
# 2. Section headers (single label cells) move from column H to column G.
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "FP"
$ws.Range("G10").Value = "FN"
$ws.Range("G19").Value = "ACC"

# ---------------------------------------------------------------------
# 3. Column sub-headers (Case 4 / Case 3 / Case 2 / Case 1), now G:J.
# ---------------------------------------------------------------------
$headerRows = @(2, 11, 20)
$headerVals = @("Case 4", "Case 3", "Case 2", "Case 1")
$dataCols = @("G", "H", "I", "J")
foreach ($hr in $headerRows) {
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Range(($dataCols[$i] + $hr)).Value = $headerVals[$i]
    }
}

# ---------------------------------------------------------------------
# 4. Raw numeric data blocks, now G3:J8 and G12:J17 (previously I3:L8 and
#    I12:L17).
# ---------------------------------------------------------------------
$block1 = @(
    @(110, 123, 134, 159),
    @(220, 244, 257, 284),
    @(213, 237, 251, 278),
    @(216, 241, 255, 282),
    @(213, 237, 250, 277),
    @(196, 218, 230, 255)
)
for ($r = 0; $r -lt 6; $r++) {
    $row = 3 + $r
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Range(($dataCols[$c] + $row)).Value = $block1[$r][$c]
    }
}

$block2 = @(
    @(81, 63, 55, 48),
    @(16, 13, 11, 10),
    @(11, 8, 7, 6),
    @(9, 7, 6, 5),
    @(15, 12, 11, 10),
    @(2, 2, 2, 2)
)
for ($r = 0; $r -lt 6; $r++) {
    $row = 12 + $r
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Range(($dataCols[$c] + $row)).Value = $block2[$r][$c]
    }
}

# ---------------------------------------------------------------------
# 5. Formula block G21:J26 (previously I21:L26), no longer shared
#    formulas - each cell gets its own plain formula.
# ---------------------------------------------------------------------
$leftCols = @("B", "C", "D", "E")
for ($r = 0; $r -lt 6; $r++) {
    $destRow = 21 + $r
    $srcRow = 12 + $r
    $dateRow = 3 + $r
    for ($c = 0; $c -lt 4; $c++) {
        $lc = $leftCols[$c]
        $dc = $dataCols[$c]
        $formula = "=100*(" + $lc + $srcRow + "+" + $lc + $destRow + ")/(" + $lc + $srcRow + "+" + $lc + $destRow + "+" + $dc + $srcRow + "+" + $dc + $dateRow + ")"
        $cell = $ws.Range(($dc + $destRow))
        $cell.Formula = $formula
        $cell.NumberFormat = "0.0"
    }
}

# ---------------------------------------------------------------------
# 6. Column widths: column A widens slightly and loses "best fit"; the
#    special-width column moves from H to G.
# ---------------------------------------------------------------------
$ws.Columns.Item("H").ColumnWidth = $ws.Columns.Item("A").ColumnWidth
$ws.Columns.Item("A").ColumnWidth = 12
$ws.Columns.Item("G").ColumnWidth = $ws.Columns.Item("H").ColumnWidth
$ws.Columns.Item("H").ColumnWidth = 8.43

# ---------------------------------------------------------------------
# 7. Reset the custom row heights on rows 24-26 back to the sheet
#    default, and drop the now-empty, custom-height spacer rows
#    28/30/31 entirely.
# ---------------------------------------------------------------------
$ws.Range("A24:A31").EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 8. Selection, used just for cosmetics but mirrors the authored file.
# ---------------------------------------------------------------------
$ws.Range("N14").Select()

# ---------------------------------------------------------------------
# 9. Page setup: landscape -> portrait, and a new print area.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PrintArea = '$A$1:$M$53'

# ---------------------------------------------------------------------
# 10. Reposition/resize the chart: it now starts near the top-left of
#     the data block and is slightly smaller.
# ---------------------------------------------------------------------
$chart = $ws.ChartObjects(1)
$chart.Top = $ws.Range("A28").Top
$chart.Left = $ws.Range("A28").Left
$chart.Width = 360
$chart.Height = 360
